$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A48").Copy()
$ws.Range("A49").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A49").Value = 45393
$ws.Range("B49").Value = 3.5
$ws.Range("C49").Formula = "=C48+B49"

$ws.Range("C49").Select()
